# Scheduled-runner style refresh of the FFXIV leve-profitability workbook.
# Re-pulls currentAveragePrice* / LevePrice* / LeveProfit* figures (columns H:N)
# for a handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# and writes the refreshed numbers back in place (no layout/formatting changes).

$wb = $excel.ActiveWorkbook

# Map of Sheet name -> Row number -> Column letter -> new value.
# A value of $null means the cell is cleared entirely (matches rows where the
# source feed stopped returning a HQ price so the old figure is removed, not zeroed).
$changes = @{
    "ALC" = @{
        6 = @{ "H" = 595.3333; "I" = 460.44446; "J" = 1000; "K" = 1381.33338; "L" = 3000; "M" = -1269.33338; "N" = -3224 }
        12 = @{ "H" = 46.5; "I" = 91; "J" = 2; "K" = 91; "L" = 2; "M" = 79; "N" = -342 }
        38 = @{ "H" = 19; "I" = 19; "K" = 57; "M" = 315 }
        58 = @{ "H" = 4449.5713; "I" = 242; "J" = 8657.143; "K" = 726; "L" = 25971.429; "M" = -576; "N" = -26271.429 }
        87 = @{ "H" = 0; "J" = 0; "L" = 0; "N" = $null }
        90 = @{ "H" = 0; "J" = 0; "L" = 0; "N" = $null }
        137 = @{ "H" = 50001920; "I" = 58824788; "J" = 5667.6665; "K" = 176474364; "L" = 17002.9995; "M" = -176471814; "N" = -22102.9995 }
        138 = @{ "H" = 8407427; "I" = 5956859.5; "J" = 8932549; "K" = 17870578.5; "L" = 26797647; "M" = -17865438.5; "N" = -26807927 }
    }
    "ARM" = @{
        32 = @{ "H" = 42051.125; "I" = 9613.807000000001; "J" = 153779.67; "K" = 9613.807000000001; "L" = 153779.67; "M" = -9326.807000000001; "N" = -154353.67 }
        45 = @{ "H" = 1427.5217; "I" = 1430.95; "K" = 1430.95; "M" = -1053.95 }
        110 = @{ "H" = 111111920; "I" = 250000660; "K" = 250000660; "M" = -249998615 }
        132 = @{ "H" = 4513.1763; "I" = 4285.5386; "J" = 5253; "K" = 12856.6158; "L" = 15759; "M" = -10326.6158; "N" = -20819 }
    }
    "BSM" = @{
        105 = @{ "H" = 2804.1365; "I" = 2705.516; "J" = 3039.3076; "K" = 2705.516; "L" = 3039.3076; "M" = -958.5160000000001; "N" = -6533.3076 }
    }
    "CRP" = @{
        31 = @{ "H" = 6025.173; "I" = 3549.818; "J" = 6689.2925; "K" = 3549.818; "L" = 6689.2925; "M" = -3254.818; "N" = -7279.2925 }
        34 = @{ "H" = 6025.173; "I" = 3549.818; "J" = 6689.2925; "K" = 3549.818; "L" = 6689.2925; "M" = -3347.818; "N" = -7093.2925 }
        134 = @{ "H" = 62503750; "I" = 500002500; "J" = 22731136; "K" = 1500007500; "L" = 68193408; "M" = -1500004965; "N" = -68198478 }
    }
    "CUL" = @{
        64 = @{ "H" = 3768.3809; "I" = 1955.3334; "J" = 4493.6; "K" = 5866.0002; "L" = 13480.8; "M" = -5596.0002; "N" = -14020.8 }
        67 = @{ "H" = 3768.3809; "I" = 1955.3334; "J" = 4493.6; "K" = 5866.0002; "L" = 13480.8; "M" = -4930.0002; "N" = -15352.8 }
        130 = @{ "H" = 915; "J" = 1300; "L" = 3900; "N" = -13940 }
    }
    "GSM" = @{
        122 = @{ "H" = 2761.5; "I" = 3032.524; "J" = 1948.4286; "K" = 9097.572; "L" = 5845.2858; "M" = -6647.572; "N" = -10745.2858 }
        132 = @{ "H" = 7146.2856; "I" = 7670.6665; "J" = 4000; "K" = 23011.9995; "L" = 12000; "M" = -20481.9995; "N" = -17060 }
    }
    "LTW" = @{
        40 = @{ "H" = 3050.2778; "I" = 1600; "J" = 3464.6428; "K" = 1600; "L" = 3464.6428; "M" = -1464; "N" = -3736.6428 }
        68 = @{ "H" = 1838.7; "I" = 1783.8572; "J" = 1966.6666; "K" = 1783.8572; "L" = 1966.6666; "M" = -1034.8572; "N" = -3464.6666 }
        71 = @{ "H" = 1838.7; "I" = 1783.8572; "J" = 1966.6666; "K" = 8919.286; "L" = 9833.333000000001; "M" = -5175.286; "N" = -17321.333 }
        122 = @{ "H" = 4073.3333; "J" = 4073.3333; "L" = 12219.9999; "N" = -17119.9999 }
    }
    "WVR" = @{
        132 = @{ "H" = 2573.3684; "I" = 2314.9148; "J" = 3788.1; "K" = 6944.7444; "L" = 11364.3; "M" = -4414.7444; "N" = -16424.3 }
    }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $changes[$sheetName]
    foreach ($rowNum in $rowsForSheet.Keys) {
        $colsForRow = $rowsForSheet[$rowNum]
        foreach ($colLetter in $colsForRow.Keys) {
            $newValue = $colsForRow[$colLetter]
            $addr = "$colLetter$rowNum"
            $ws.Range($addr).Value = $newValue
        }
    }
}

Write-Host "Refreshed profit figures across $($changes.Count) sheet(s)."
